# "Reformatting CYRS,HSI,SRS Documents & Updatting RTM"
#
# The J1 header cell previously held the shared string "Status". The
# column is being dropped from the header row (the RTM no longer tracks a
# separate "Status" column) while keeping the header's visual formatting
# footprint: a centered, bold, solid-filled cell — just switched to a
# light (white/background1) fill instead of the black header fill used by
# the rest of row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "Status" text from J1 so the shared-string table no longer
# carries that entry (Excel will compact sharedStrings.xml on save).
$hdr = $ws.Range("J1")
$hdr.ClearContents()

# Re-style J1: solid fill using the light/background-1 theme color, still
# centered like the other header cells (font stays bold/white, inherited
# from the existing header formatting already on this cell).
$hdr.Interior.Pattern = 1            # xlSolid
$hdr.Interior.ThemeColor = 2         # xlThemeColorLight1 (theme index 0 / "Background 1")
$hdr.HorizontalAlignment = -4108     # xlCenter

# Update the active selection on the sheet to reflect where the author
# left off editing.
$ws.Range("J3").Select()
